$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to stay text so numeric-looking strings
# like "599.30" or "5.86" are not coerced into floating point numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '65.658.84'
$ws.Range("E2").Value = '  -0.31%  '
$ws.Range("D3").Value = '2.667.78'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '599.30'
$ws.Range("E5").Value = '  -1.45%  '
$ws.Range("D6").Value = '156.59'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +4.88%  '
$ws.Range("D9").Value = '0.130'
$ws.Range("E9").Value = '  +4.78%  '
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("D11").Value = '5.86'
$ws.Range("E11").Value = '  -2.50%  '
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = '29.30'
$ws.Range("E13").Value = '  -3.47%  '
$ws.Range("D14").Value = '0.0000197'
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("D15").Value = '3.146.53'
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("D16").Value = '65.482.57'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '2.659.84'
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("D18").Value = '12.81'
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("D19").Value = '4.78'
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").Value = '7.56'
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").Value = '351.44'
$ws.Range("E21").Value = '  -2.19%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '69.50'
$ws.Range("E23").Value = '  -1.67%  '
$ws.Range("D24").Value = '0.0000112'
$ws.Range("E24").Value = '  +4.64%  '
$ws.Range("D25").Value = '9.59'
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D26").Value = '1.64'
$ws.Range("E26").Value = '  -2.00%  '
$ws.Range("E27").Value = '  -2.63%  '
$ws.Range("E28").Value = '  -5.76%  '
$ws.Range("D29").Value = '8.01'
$ws.Range("E29").Value = '  -4.84%  '
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("D32").Value = '525.32'
$ws.Range("E32").Value = '  -2.91%  '
$ws.Range("E33").Value = '  -2.45%  '
$ws.Range("D34").Value = '6.45'
$ws.Range("E34").Value = '  -3.30%  '
$ws.Range("D35").Value = '5.48'
$ws.Range("E35").Value = '  +1.58%  '
$ws.Range("D36").Value = '0.422'
$ws.Range("E36").Value = '  -2.45%  '
$ws.Range("D37").Value = '20.53'
$ws.Range("E37").Value = '  -1.48%  '
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = '158.01'
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("E40").Value = '  -2.99%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '163.22'
$ws.Range("E42").Value = '  -4.92%  '
$ws.Range("D43").Value = '4.12'
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("D44").Value = '2.31'
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("D45").Value = '0.0608'
$ws.Range("E45").Value = '  -1.06%  '
$ws.Range("D46").Value = '22.74'
$ws.Range("E46").Value = '  -3.38%  '
$ws.Range("D47").Value = '0.639'
$ws.Range("E47").Value = '  -2.73%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0257'
$ws.Range("E48").Value = '  -3.28%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0263'
$ws.Range("E49").Value = '  +14.16%  '
$ws.Range("D50").Value = '0.100'
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("D51").Value = '20.06'
$ws.Range("E51").Value = '  -4.54%  '

# Restore the default cell style on column D so no stray number format
# style sticks around on cells that did not have one originally.
$priceRange.Style = "Normal"
